$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 75
$ws.Range("I5").Value = 70.833336
$ws.Range("K5").Value = 70.833336
$ws.Range("M5").Value = 44.166664
$ws.Range("H19").Value = 1965.8
$ws.Range("I19").Value = 1958.5
$ws.Range("J19").Value = 1995
$ws.Range("K19").Value = 1958.5
$ws.Range("L19").Value = 1995
$ws.Range("M19").Value = -1783.5
$ws.Range("N19").Value = -2345
$ws.Range("H32").Value = 1498.5
$ws.Range("J32").Value = 1498.5
$ws.Range("L32").Value = 1498.5
$ws.Range("N32").Value = -2150.5
$ws.Range("H40").Value = 2000
$ws.Range("J40").Value = 2000
$ws.Range("L40").Value = 2000
$ws.Range("N40").Value = -2350
$ws.Range("H98").Value = 7639.316
$ws.Range("I98").Value = 2164.7
$ws.Range("K98").Value = 2164.7
$ws.Range("M98").Value = -666.6999999999998
$ws.Range("H122").Value = 7639.316
$ws.Range("I122").Value = 2164.7
$ws.Range("K122").Value = 6494.099999999999
$ws.Range("M122").Value = -4044.099999999999
$ws.Range("H125").Value = 1691.3334
$ws.Range("J125").Value = 1871.5
$ws.Range("L125").Value = 16843.5
$ws.Range("N125").Value = -21763.5
$ws.Range("H137").Value = 11217.917
$ws.Range("J137").Value = 13699
$ws.Range("L137").Value = 41097
$ws.Range("N137").Value = -46197

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 9995
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H45").Value = 9166
$ws.Range("I45").Value = 4749
$ws.Range("K45").Value = 4749
$ws.Range("M45").Value = -4372
$ws.Range("H88").Value = 1314.625
$ws.Range("I88").Value = 924.3333
$ws.Range("K88").Value = 924.3333
$ws.Range("M88").Value = -518.3333
$ws.Range("H91").Value = 1314.625
$ws.Range("I91").Value = 924.3333
$ws.Range("K91").Value = 924.3333
$ws.Range("M91").Value = 479.6667
$ws.Range("H132").Value = 11999.8
$ws.Range("I132").Value = 5999.5
$ws.Range("J132").Value = 16000
$ws.Range("K132").Value = 17998.5
$ws.Range("L132").Value = 48000
$ws.Range("M132").Value = -15468.5
$ws.Range("N132").Value = -53060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2170.3845
$ws.Range("I20").Value = 2110.5454
$ws.Range("K20").Value = 2110.5454
$ws.Range("M20").Value = -1863.5454
$ws.Range("H86").Value = 1978.4445
$ws.Range("I86").Value = 1971.2858
$ws.Range("J86").Value = 2003.5
$ws.Range("K86").Value = 1971.2858
$ws.Range("L86").Value = 2003.5
$ws.Range("M86").Value = -848.2858000000001
$ws.Range("N86").Value = -4249.5
$ws.Range("H89").Value = 1978.4445
$ws.Range("I89").Value = 1971.2858
$ws.Range("J89").Value = 2003.5
$ws.Range("K89").Value = 9856.429
$ws.Range("L89").Value = 10017.5
$ws.Range("M89").Value = -4240.429
$ws.Range("N89").Value = -21249.5
$ws.Range("H99").Value = 3167
$ws.Range("I99").Value = 3167
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3167
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -1669
$ws.Range("H107").Value = 2827.5557
$ws.Range("I107").Value = 778.2857
$ws.Range("K107").Value = 778.2857
$ws.Range("M107").Value = 1141.7143
$ws.Range("H134").Value = 9803.167
$ws.Range("I134").Value = 4704.75
$ws.Range("J134").Value = 20000
$ws.Range("K134").Value = 14114.25
$ws.Range("L134").Value = 60000
$ws.Range("M134").Value = -11579.25
$ws.Range("N134").Value = -65070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 79.95238
$ws.Range("I7").Value = 36.142857
$ws.Range("J7").Value = 101.85714
$ws.Range("K7").Value = 36.142857
$ws.Range("L7").Value = 101.85714
$ws.Range("M7").Value = 76.85714300000001
$ws.Range("N7").Value = -327.85714
$ws.Range("H22").Value = 322.6
$ws.Range("I22").Value = 232.28572
$ws.Range("K22").Value = 232.28572
$ws.Range("M22").Value = 117.71428
$ws.Range("H31").Value = 6699.5
$ws.Range("I31").Value = 5000
$ws.Range("J31").Value = 8399
$ws.Range("K31").Value = 5000
$ws.Range("L31").Value = 8399
$ws.Range("M31").Value = -4705
$ws.Range("N31").Value = -8989
$ws.Range("H34").Value = 6699.5
$ws.Range("I34").Value = 5000
$ws.Range("J34").Value = 8399
$ws.Range("K34").Value = 5000
$ws.Range("L34").Value = 8399
$ws.Range("M34").Value = -4798
$ws.Range("N34").Value = -8803
$ws.Range("H132").Value = 5239.4
$ws.Range("I132").Value = 2770.5715
$ws.Range("J132").Value = 11000
$ws.Range("K132").Value = 8311.7145
$ws.Range("L132").Value = 33000
$ws.Range("M132").Value = -5781.7145
$ws.Range("N132").Value = -38060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 166835040
$ws.Range("I4").Value = 202049.4
$ws.Range("K4").Value = 606148.2
$ws.Range("M4").Value = -606036.2
$ws.Range("H23").Value = 111.22222
$ws.Range("I23").Value = 152
$ws.Range("J23").Value = 78.6
$ws.Range("K23").Value = 456
$ws.Range("L23").Value = 235.8
$ws.Range("M23").Value = -221
$ws.Range("N23").Value = -705.8
$ws.Range("H40").Value = 119
$ws.Range("I40").Value = 65
$ws.Range("J40").Value = 200
$ws.Range("K40").Value = 260
$ws.Range("L40").Value = 800
$ws.Range("M40").Value = -191
$ws.Range("N40").Value = -938
$ws.Range("H137").Value = 500
$ws.Range("I137").Value = 500
$ws.Range("K137").Value = 1500
$ws.Range("M137").Value = 3600

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").ClearContents()
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = 0
$ws.Range("H41").Value = 1400
$ws.Range("I41").Value = 1400
$ws.Range("K41").Value = 1400
$ws.Range("M41").Value = -1045
$ws.Range("H122").Value = 741.4286
$ws.Range("I122").Value = 719
$ws.Range("J122").Value = 797.5
$ws.Range("K122").Value = 2157
$ws.Range("L122").Value = 2392.5
$ws.Range("M122").Value = 293
$ws.Range("N122").Value = -7292.5
$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 7656.7144
$ws.Range("I132").Value = 4149.25
$ws.Range("J132").Value = 12333.333
$ws.Range("K132").Value = 12447.75
$ws.Range("L132").Value = 36999.999
$ws.Range("M132").Value = -9917.75
$ws.Range("N132").Value = -42059.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H7").Value = 5499
$ws.Range("I7").Value = 5499
$ws.Range("K7").Value = 5499
$ws.Range("M7").Value = -5387
$ws.Range("H16").Value = 1599.6
$ws.Range("I16").Value = 1332.6666
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1332.6666
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1162.6666
$ws.Range("N16").Value = -2340
$ws.Range("H22").Value = 1450.9333
$ws.Range("I22").Value = 789.1429
$ws.Range("J22").Value = 2030
$ws.Range("K22").Value = 789.1429
$ws.Range("L22").Value = 2030
$ws.Range("M22").Value = -494.1429000000001
$ws.Range("N22").Value = -2620
$ws.Range("H27").Value = 1450.9333
$ws.Range("I27").Value = 789.1429
$ws.Range("J27").Value = 2030
$ws.Range("K27").Value = 789.1429
$ws.Range("L27").Value = 2030
$ws.Range("M27").Value = -682.1429
$ws.Range("N27").Value = -2244
$ws.Range("H40").Value = 1987.75
$ws.Range("I40").Value = 2128.1428
$ws.Range("J40").Value = 1005
$ws.Range("K40").Value = 2128.1428
$ws.Range("L40").Value = 1005
$ws.Range("M40").Value = -1992.1428
$ws.Range("N40").Value = -1277
$ws.Range("H46").Value = 5002
$ws.Range("J46").Value = 5002
$ws.Range("L46").Value = 5002
$ws.Range("N46").Value = -5378
$ws.Range("H61").Value = 3664.3333
$ws.Range("I61").Value = 3664.3333
$ws.Range("K61").Value = 3664.3333
$ws.Range("M61").Value = -3462.3333
$ws.Range("H113").Value = 3664.3333
$ws.Range("I113").Value = 3664.3333
$ws.Range("K113").Value = 3664.3333
$ws.Range("M113").Value = -1494.3333
$ws.Range("H122").Value = 3539
$ws.Range("I122").Value = 3388.5
$ws.Range("K122").Value = 10165.5
$ws.Range("M122").Value = -7715.5
$ws.Range("H126").Value = 5499
$ws.Range("I126").Value = 5499
$ws.Range("K126").Value = 16497
$ws.Range("M126").Value = -14027

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 19985
$ws.Range("I58").Value = 19985
$ws.Range("K58").Value = 19985
$ws.Range("M58").Value = -19677
$ws.Range("H80").Value = 69000
$ws.Range("J80").Value = 69000
$ws.Range("L80").Value = 69000
$ws.Range("N80").Value = -70996
$ws.Range("H83").Value = 69000
$ws.Range("J83").Value = 69000
$ws.Range("L83").Value = 207000
$ws.Range("N83").Value = -216984
$ws.Range("H122").Value = 2366.1667
$ws.Range("I122").Value = 2439.4
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 7318.200000000001
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -4868.200000000001
$ws.Range("N122").Value = -10900
$ws.Range("H126").Value = 1899.1666
$ws.Range("I126").Value = 2039
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 6117
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = -3647
$ws.Range("N126").Value = -8540
$ws.Range("H132").Value = 7593.909
$ws.Range("I132").Value = 4797.857
$ws.Range("J132").Value = 12487
$ws.Range("K132").Value = 14393.571
$ws.Range("L132").Value = 37461
$ws.Range("M132").Value = -11863.571
$ws.Range("N132").Value = -42521
